# AGP - Checklist Verificacao de Projeto
# Manutenção dos documentos do projeto para fase de elaboração:
# fill in the verification results on the "Ver-Elaboração1" checklist.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ver-Elaboração1")

# Items marked "Sim" (evidence found).
$ws.Range("D6").Value = "Sim"
$ws.Range("D8").Value = "Sim"
$ws.Range("D10:D13").Value = "Sim"
$ws.Range("D17").Value = "Sim"
$ws.Range("D19").Value = "Sim"
$ws.Range("D21").Value = "Sim"
$ws.Range("D23").Value = "Sim"
$ws.Range("D28:D33").Value = "Sim"
$ws.Range("D35:D40").Value = "Sim"
$ws.Range("D42:D44").Value = "Sim"

# Items marked "Não" (no evidence).
$ws.Range("D25:D26").Value = "Não"

# Items marked "NA" (not applicable).
$ws.Range("D15").Value = "NA"
$ws.Range("D46:D48").Value = "NA"

# Recalculate so every dependent formula (counts, percentages, the
# "Indicadores"/"Detalhado" summary sheets and their charts) picks up
# the new inputs.
$excel.CalculateFullRebuild()

# Leave the workbook focused on this sheet / the last-edited cell, as it
# was when the author saved.
$ws.Activate()
$ws.Range("D48").Select()
